# Automatische test-sync: 2025-06-30 20:12:50
#
# Adds the new "Wanneer krijg ik mijn offerte?" test-mail row to the
# "Logs" sheet (row 16), extends the conditional-formatting ranges that
# covered rows 2:15 to now cover rows 2:16, and updates the "Dashboard"
# sheet summary rows 4-6 to reflect the new category counts / ordering.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append row 16 with the new test-mail entry.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(16, 1).Value = "Wanneer krijg ik mijn offerte?"
$logs.Cells.Item(16, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(16, 3).Value = "Testmail #16: Wanneer krijg ik mijn offerte?"
$logs.Cells.Item(16, 4).Value = "Offerte / Prijsaanvraag"
$logs.Cells.Item(16, 5).Value = "Geachte klant,`nDank u voor uw e-mail. Uw offerte zal naar verwachting binnen 24 uur worden verstuurd. Mocht u deze niet op tijd ontvangen, neem dan gerust contact met ons op.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item(16, 6).Value = "2025-06-30 20:12:43"
$logs.Cells.Item(16, 7).Value = "Ja"
$logs.Cells.Item(16, 8).Value = "Nee"
$logs.Cells.Item(16, 9).Value = "Ja"
$logs.Cells.Item(16, 10).Value = "Ja"

# ---------------------------------------------------------------------
# 2) Extend the conditional-formatting "appliesTo" ranges from row 15
#    to row 16 for columns D, G, H, I, J.
# ---------------------------------------------------------------------
$ctRanges = @("D", "G", "H", "I", "J")
foreach ($col in $ctRanges) {
    $fcs = $logs.Range("$col" + "2:" + "$col" + "15").FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($logs.Range("$col" + "2:" + "$col" + "16"))
    }
}

# ---------------------------------------------------------------------
# 3) Dashboard sheet: the new mail bumps "Offerte / Prijsaanvraag" to 2
#    occurrences, which now outranks "Bestelling / Levering" and
#    "Openingstijden / Locatie" (both still at 1), re-sorting rows 4-6.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Offerte / Prijsaanvraag"
$dash.Cells.Item(4, 2).Value = 2

$dash.Cells.Item(5, 1).Value = "Bestelling / Levering"
$dash.Cells.Item(5, 2).Value = 1

$dash.Cells.Item(6, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(6, 2).Value = 1
